$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.012.75'
$ws.Range('E2').Value = '  -3.14%  '
$ws.Range('D3').Value = '2.361.96'
$ws.Range('E3').Value = '  -4.07%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = "'" + '499.73'
$ws.Range('E5').Value = '  -2.28%  '
$ws.Range('D6').Value = "'" + '129.21'
$ws.Range('E6').Value = '  -3.72%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('E8').Value = '  -2.48%  '
$ws.Range('D9').Value = '2.364.83'
$ws.Range('E9').Value = '  -3.94%  '
$ws.Range('D10').Value = "'" + '0.0980'
$ws.Range('E10').Value = '  -0.36%  '
$ws.Range('E11').Value = '  +0.27%  '
$ws.Range('D12').Value = "'" + '4.78'
$ws.Range('E12').Value = '  +2.74%  '
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').Value = '2.781.90'
$ws.Range('E14').Value = '  -3.94%  '
$ws.Range('D15').Value = '55.992.47'
$ws.Range('E15').Value = '  -3.14%  '
$ws.Range('D16').Value = "'" + '21.39'
$ws.Range('E16').Value = '  -2.87%  '
$ws.Range('E17').Value = '  -2.17%  '
$ws.Range('D18').Value = '2.311.63'
$ws.Range('E18').Value = '  -6.97%  '
$ws.Range('D19').Value = "'" + '10.00'
$ws.Range('E19').Value = '  -3.59%  '
$ws.Range('D20').Value = "'" + '4.03'
$ws.Range('E20').Value = '  -3.40%  '
$ws.Range('D21').Value = "'" + '306.78'
$ws.Range('E21').Value = '  -2.79%  '
$ws.Range('E22').Value = '  -3.12%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = "'" + '65.35'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').Value = "'" + '0.998'
$ws.Range('E25').Value = '  -0.12%  '
$ws.Range('E26').Value = '  -3.88%  '
$ws.Range('E27').Value = '  -6.14%  '
$ws.Range('D28').Value = "'" + '7.22'
$ws.Range('E28').Value = '  -5.07%  '
$ws.Range('D29').Value = "'" + '171.19'
$ws.Range('E29').Value = '  -0.89%  '
$ws.Range('D30').Value = '0.0₃0709'
$ws.Range('E30').Value = '  -3.63%  '
$ws.Range('E31').Value = '  -3.64%  '
$ws.Range('E32').Value = '  +0.04%  '
$ws.Range('D33').Value = "'" + '5.75'
$ws.Range('D34').Value = "'" + '0.998'
$ws.Range('E34').Value = '  +0.12%  '
$ws.Range('E35').Value = '  -5.45%  '
$ws.Range('D36').Value = "'" + '17.58'
$ws.Range('E36').Value = '  -2.90%  '
$ws.Range('E37').Value = '  -5.85%  '
$ws.Range('E38').Value = '  -4.15%  '
$ws.Range('D39').Value = "'" + '36.10'
$ws.Range('E39').Value = '  -1.79%  '
$ws.Range('D40').Value = "'" + '0.789'
$ws.Range('E40').Value = '  -2.90%  '
$ws.Range('E41').Value = '  -5.93%  '
$ws.Range('E42').Value = '  -2.10%  '
$ws.Range('D43').Value = "'" + '128.84'
$ws.Range('E43').Value = '  -6.32%  '
$ws.Range('D44').Value = "'" + '4.70'
$ws.Range('E44').Value = '  -4.46%  '
$ws.Range('D45').Value = "'" + '0.562'
$ws.Range('E45').Value = '  -2.75%  '
$ws.Range('E46').Value = '  -2.34%  '
$ws.Range('D47').Value = "'" + '239.02'
$ws.Range('E47').Value = '  -7.21%  '
$ws.Range('E48').Value = '  -2.94%  '
$ws.Range('E49').Value = '  -4.21%  '
$ws.Range('E50').Value = '  -2.38%  '
